$d = $word.ActiveDocument

# Five target paragraphs that become the new document body:
#   1) Title "Trip Planner App"  - centered, 20pt      (sz/szCs 40)
#   2) "Overview:" heading       - 14pt                (sz/szCs 28)
#   3) Body paragraph #1         - 12pt, two runs       (sz/szCs 24)
#   4) Body paragraph #2         - 12pt                (sz/szCs 24)
#   5) Trailing blank paragraph  - 14pt paragraph-mark  (sz/szCs 28)
#
# Every paragraph carries <w:bidi w:val="0"/>, matching the bidi setting
# already present on the document's sole (empty) paragraph.
$bodyXml = '<w:p><w:pPr><w:bidi w:val="0"/><w:jc w:val="center"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>Trip Planner App</w:t></w:r></w:p><w:p><w:pPr><w:bidi w:val="0"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Overview:</w:t></w:r></w:p><w:p><w:pPr><w:bidi w:val="0"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> platform with which users can plan a trip and maximize vacation time according to their preferences of attractions.</w:t></w:r></w:p><w:p><w:pPr><w:bidi w:val="0"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The platform will accept the users’ preferences and vacation time limit and will calculate the optimal itinerary for them to take.</w:t></w:r></w:p><w:p><w:pPr><w:bidi w:val="0"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p>'

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Insert the five new paragraphs at the very start of the document.
$insertionPoint = $d.Range(0, 0)
[void]$insertionPoint.InsertXML($xml)

# The original (now last) paragraph is the blank paragraph that used to be
# the document's only paragraph. Our freshly-inserted XML already supplies
# its replacement (paragraph 5 above), so remove the old leftover paragraph,
# merging its trailing paragraph mark away.
$oldLast = $d.Paragraphs($d.Paragraphs.Count)
$removeRange = $d.Range($oldLast.Range.Start - 1, $oldLast.Range.End)
$removeRange.Delete()

Write-Host "Paragraphs:" $d.Paragraphs.Count
Write-Host "Text:" $d.Content.Text
